$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-24 Monday" "2025-02-25 Tuesday"

Replace-Text "23×77=1771" "29×60=1740"
Replace-Text "95×54=5130" "91×70=6370"
Replace-Text "65×31=2015" "12×86=1032"
Replace-Text "21×35=735" "13×54=702"
Replace-Text "41×85=3485" "21×90=1890"

Replace-Text "58×51=2958" "99×73=7227"
Replace-Text "26×78=2028" "13×47=611"
Replace-Text "68×60=4080" "35×26=910"
Replace-Text "95×80=7600" "75×41=3075"
Replace-Text "40×85=3400" "58×37=2146"

Replace-Text "54×90=4860" "15×79=1185"
Replace-Text "57×28=1596" "52×57=2964"
Replace-Text "51×14=714" "65×60=3900"
Replace-Text "51×88=4488" "20×99=1980"
Replace-Text "61×83=5063" "78×88=6864"

Replace-Text "87×66=5742" "95×57=5415"
Replace-Text "95×67=6365" "12×15=180"
Replace-Text "58×97=5626" "55×90=4950"
Replace-Text "82×60=4920" "47×43=2021"
Replace-Text "96×81=7776" "33×74=2442"

Replace-Text "27×33=891" "82×89=7298"
Replace-Text "32×28=896" "43×35=1505"
Replace-Text "71×16=1136" "31×65=2015"
Replace-Text "91×57=5187" "78×49=3822"
Replace-Text "99×23=2277" "29×64=1856"
